# "added 4wk low sales check"
# Updates forecast numbers (and a couple of derived status labels) on the
# "Forecast Comparison" sheet, and refreshes the dependent roll-up figures
# on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Forecast Comparison"
# ---------------------------------------------------------------------
$fc = $wb.Worksheets.Item("Forecast Comparison")

# Week W10 (row 2)
$fc.Range("D2").Value = 30
$fc.Range("H2").Value = 2.27

# Week W11 (row 3)
$fc.Range("D3").Value = 29
$fc.Range("H3").Value = 1.31
$fc.Range("L3").Value = 1.11

# Week W12 (row 4)
$fc.Range("D4").Value = 28
$fc.Range("H4").Value = 0.32
$fc.Range("I4").Value = "High"
$fc.Range("J4").Value = "Urgent"
$fc.Range("L4").Value = 0.87

# Week W13 (row 5)
$fc.Range("D5").Value = 28
$fc.Range("H5").Value = 0
$fc.Range("L5").Value = 0.95

# Week W14 (row 6)
$fc.Range("D6").Value = 27
$fc.Range("L6").Value = 0.95

# Week W15 (row 7)
$fc.Range("D7").Value = 26
$fc.Range("L7").Value = 1.14

# Week W16 (row 8)
$fc.Range("D8").Value = 26
$fc.Range("L8").Value = 1.04

# Week W17 (row 9)
$fc.Range("L9").Value = 0.83

# Week W18 (row 10)
$fc.Range("D10").Value = 24
$fc.Range("L10").Value = 0.94

# Week W19 (row 11)
$fc.Range("D11").Value = 24
$fc.Range("L11").Value = 1.03

# Week W20 (row 12)
$fc.Range("L12").Value = 1.08

# Week W21 (row 13)
$fc.Range("D13").Value = 22
$fc.Range("L13").Value = 1.18

# Week W22 (row 14)
$fc.Range("D14").Value = 21
$fc.Range("L14").Value = 0.85

# Week W23 (row 15)
$fc.Range("D15").Value = 20
$fc.Range("L15").Value = 1.06

# Week W24 (row 16)
$fc.Range("L16").Value = 1.09

# Week W25 (row 17)
$fc.Range("D17").Value = 18
$fc.Range("L17").Value = 1.13

# ---------------------------------------------------------------------
# Sheet 2: "Summary" - recomputed roll-up values (stored as text, as in
# the original workbook). A leading apostrophe keeps them as text cells
# instead of being auto-coerced to numbers.
# ---------------------------------------------------------------------
$sm = $wb.Worksheets.Item("Summary")

$sm.Range("B9").Value  = "'392"   # Total Forecast (16 Weeks)
$sm.Range("B10").Value = "'219"   # Total Forecast (8 Weeks)
$sm.Range("B11").Value = "'115"   # Total Forecast (4 Weeks)
$sm.Range("B12").Value = "'30"    # Max Forecast
$sm.Range("B14").Value = "'18"    # Min Forecast
